$wb = $excel.ActiveWorkbook
$loginSheet = $wb.Worksheets.Item("Login")

# --- Insert new sheet "Login_DataProvider" right after "Login" ---
$newSheet = $wb.Worksheets.Add($null, $loginSheet)
$newSheet.Name = "Login_DataProvider"

# Re-fetch the Customer sheet reference now that the sheet order changed
$customerSheet = $wb.Worksheets.Item("Customer")

# Copy header row (values + style) from Login!A1:B1
$loginSheet.Range("A1:B1").Copy($newSheet.Range("A1:B1"))

# Add the hyperlinks first (so the cell style gets overwritten by the
# subsequent style copy instead of the default Hyperlink style)
$newSheet.Hyperlinks.Add($newSheet.Range("A2"), "mailto:admin@example.com") | Out-Null
$newSheet.Hyperlinks.Add($newSheet.Range("A3"), "mailto:user@example.com") | Out-Null
$newSheet.Hyperlinks.Add($newSheet.Range("A4"), "mailto:admin2@example.com") | Out-Null
$newSheet.Hyperlinks.Add($newSheet.Range("A5"), "mailto:admin3@example.com") | Out-Null

# Copy the style used for the email column (Login!A2) onto the new rows
$loginSheet.Range("A2").Copy($newSheet.Range("A2:A5"))

# Fill in the data (this also re-applies the quote-prefixed-text style
# used for the password column when values start with an apostrophe)
$newSheet.Range("A2").Value = "admin@example.com"
$newSheet.Range("B2").Value = "'123456"

$newSheet.Range("A3").Value = "user@example.com"
$newSheet.Range("B3").Value = "'123"

$newSheet.Range("A4").Value = "admin2@example.com"
$newSheet.Range("B4").Value = "'A12345"

$newSheet.Range("A5").Value = "admin3@example.com"
$newSheet.Range("B5").Value = "'A12345"

# Match the page setup (paper size / orientation) used for the new sheet
$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

# --- Minor formatting touch-up on the Login sheet (matches a style
# clean-up that happens when the workbook is re-saved) ---
$loginSheet.Range("C2").Interior.Pattern = -4142
$loginSheet.Range("C3").Interior.Pattern = -4142

# --- Update selections on the various sheets ---
$loginSheet.Range("C5").Select()
$customerSheet.Range("C8").Select()
$newSheet.Range("C15").Select()

Write-Host "done"
